$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = '29.101.70'
$ws.Cells.Item(2, 5).Value = '  -1.14%  '
$ws.Cells.Item(3, 4).Value = '1.831.27'
$ws.Cells.Item(3, 5).Value = '  -1.29%  '
$ws.Cells.Item(4, 5).Value = '  -0.09%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '239.36'
$ws.Cells.Item(5, 5).Value = '  -2.25%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.6629'
$ws.Cells.Item(6, 5).Value = '  -4.63%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '1.0000'
$ws.Cells.Item(7, 5).Value = '  -0.05%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.2943'
$ws.Cells.Item(8, 5).Value = '  -3.97%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.07313'
$ws.Cells.Item(9, 5).Value = '  -4.72%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '22.67'
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.07660'
$ws.Cells.Item(11, 5).Value = '  -1.47%  '
$ws.Cells.Item(12, 4).Value = '1.834.33'
$ws.Cells.Item(12, 5).Value = '  -1.15%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '5.018'
$ws.Cells.Item(13, 5).Value = '  -2.49%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '0.6736'
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '85.79'
$ws.Cells.Item(15, 5).Value = '  -5.71%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '6.117'
$ws.Cells.Item(16, 5).Value = '  -2.87%  '
$ws.Cells.Item(17, 4).Value = '29.086.33'
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.000008208'
$ws.Cells.Item(18, 5).Value = '  -1.27%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '226.99'
$ws.Cells.Item(19, 5).Value = '  -4.62%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '12.47'
$ws.Cells.Item(20, 5).Value = '  -1.89%  '
$ws.Cells.Item(21, 5).Value = '  -0.08%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '7.261'
$ws.Cells.Item(22, 5).Value = '  -4.62%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '0.9998'
$ws.Cells.Item(23, 5).Value = '  -0.02%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '160.55'
$ws.Cells.Item(24, 5).Value = '  +0.21%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.1415'
$ws.Cells.Item(25, 5).Value = '  -5.14%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '8.646'
$ws.Cells.Item(26, 5).Value = '  -2.56%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '17.92'
$ws.Cells.Item(27, 5).Value = '  -1.76%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '1.499'
$ws.Cells.Item(28, 5).Value = '  -2.03%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '4.222'
$ws.Cells.Item(29, 5).Value = '  -0.40%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '4.086'
$ws.Cells.Item(30, 5).Value = '  -1.53%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '1.198'
$ws.Cells.Item(31, 5).Value = '  -1.42%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '0.05308'
$ws.Cells.Item(32, 5).Value = '  +4.19%  '
$ws.Cells.Item(33, 5).Value = '  -1.13%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.7446'
$ws.Cells.Item(34, 5).Value = '  -3.35%  '
$ws.Cells.Item(35, 5).Value = '  -1.94%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.678'
$ws.Cells.Item(36, 5).Value = '  -0.05%  '
$ws.Cells.Item(37, 4).Value = '1.306.28'
$ws.Cells.Item(37, 5).Value = '  -1.95%  '
$ws.Cells.Item(38, 5).Value = '  -3.70%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '2.707'
$ws.Cells.Item(39, 5).Value = '  -0.56%  '
$ws.Cells.Item(40, 5).Value = '  -2.67%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '6.029'
$ws.Cells.Item(41, 5).Value = '  +3.65%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.9981'
$ws.Cells.Item(42, 5).Value = '  -0.25%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '103.31'
$ws.Cells.Item(44, 4).Value = '1.982.10'
$ws.Cells.Item(44, 5).Value = '  -0.79%  '
$ws.Cells.Item(46, 5).Value = '  -3.20%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '63.46'
$ws.Cells.Item(47, 5).Value = '  +0.42%  '
$ws.Cells.Item(48, 5).Value = '  -1.53%  '
$ws.Cells.Item(49, 2).Value = 'XinFinNetwork'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.07593'
$ws.Cells.Item(49, 5).Value = '  +13.21%  '
$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '9.273'
$ws.Cells.Item(50, 5).Value = '  -5.82%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.05917'
$ws.Cells.Item(51, 5).Value = '  -0.29%  '
